$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1): add P1 = 14, Q1 = 15, copying the same
# formatting (style) used by the existing header cells (e.g. O1).
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Data rows 2-25: swap the I/K and M/O column values, and add the new
# P/Q columns (value 2 for every row).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2  (was 1)
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1  (was 2)
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2  (was 1)
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1  (was 2)
    $ws.Cells.Item($r, 16).Value = 2  # P (new) = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q (new) = 2
}
